$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the G column values for the two remaining data rows
$ws.Range("G4").Value = 146
$ws.Range("G5").Value = 127

# Add new column H with 2020 data
$ws.Range("H3").Value = 2020
$ws.Range("H4").Value = 158
$ws.Range("H5").Value = 397

# Extend the separator row border into the new column H
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)

# Remove rows 6 and 7 (Похищенные/Abducted/Уурдалган and Утерянные/Lost/Жоготулган)
$ws.Rows("6:7").Delete()
